# Add an "Done" column (D) to the feature tracking sheet and record a
# completion date for the "Commit code to GitHub" feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in D1, matching the style used by the other header cells (C1).
$ws.Range("D1").Value = "Done"
$ws.Range("D1").Style = $ws.Range("C1").Style

# Mark F3 ("Commit code to GitHub", row 4) as done on 2015-05-29.
$ws.Range("D4").Value = 42153
$ws.Range("D4").NumberFormat = "mm-dd-yy"

# Move the active selection to D8, as recorded in the saved workbook.
$ws.Range("D8").Select() | Out-Null
